$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data stores runs/balls/fours/sixes as text-typed numeric
# strings (Excel shows the "number stored as text" indicator). Setting
# NumberFormat to "@" before assigning the value keeps each cell text-typed
# instead of Excel auto-converting the literal to a numeric value.
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "34"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "12"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1"
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "41"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "25"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "2"
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "0"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "2"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0"
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = "0"
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "11"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "14"
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = "0"
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "18"
$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value = "1"
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "9"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "4"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "2"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "13"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "7"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "60"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "3"
$ws.Range("F10").NumberFormat = "@"
$ws.Range("F10").Value = "5"
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "47"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "3"
$ws.Range("F11").NumberFormat = "@"
$ws.Range("F11").Value = "4"
